$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.580.53"
$ws.Range("E2").Value = "  +1.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.29"
$ws.Range("E3").Value = "  +1.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.44"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5316"
$ws.Range("E7").Value = "  -2.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3960"
$ws.Range("E8").Value = "  +4.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07755"
$ws.Range("E9").Value = "  +4.16%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.02"
$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.117"
$ws.Range("E11").Value = "  +2.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.16"
$ws.Range("E12").Value = "  +3.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.303"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.579"
$ws.Range("E14").Value = "  +3.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.002"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.824.28"
$ws.Range("E16").Value = "  +1.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.15"
$ws.Range("E17").Value = "  +3.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001087"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06615"
$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.76"
$ws.Range("E20").Value = "  +1.92%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.084"
$ws.Range("E22").Value = "  +2.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.588.93"
$ws.Range("E23").Value = "  +1.62%  "

$ws.Range("E24").Value = "  +0.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.235"
$ws.Range("E25").Value = "  +6.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.73"
$ws.Range("E26").Value = "  +1.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.14"
$ws.Range("E27").Value = "  +1.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.035.37"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.420"
$ws.Range("E29").Value = "  +4.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.63"
$ws.Range("E30").Value = "  +3.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.144"
$ws.Range("E31").Value = "  +2.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1126"
$ws.Range("E32").Value = "  +1.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.734"
$ws.Range("E33").Value = "  +3.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.655"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07334"
$ws.Range("E35").Value = "  +5.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2263"
$ws.Range("E36").Value = "  +1.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02352"
$ws.Range("E37").Value = "  +2.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.895"
$ws.Range("E38").Value = "  +5.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.198"
$ws.Range("E39").Value = "  +2.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.38"
$ws.Range("E40").Value = "  +1.97%  "

$ws.Range("E41").Value = "  +2.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.195"
$ws.Range("E42").Value = "  +2.09%  "

$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.394"
$ws.Range("E44").Value = "  -2.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.47"
$ws.Range("E45").Value = "  +1.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5940"
$ws.Range("E46").Value = "  +3.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.720"
$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.39"
$ws.Range("E48").Value = "  +0.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.001"
$ws.Range("E49").Value = "  +4.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.192"
$ws.Range("E50").Value = "  +0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06952"
$ws.Range("E51").Value = "  +1.96%  "
